$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "summ30005972"
$wb.Worksheets.Item(2).Name = "summ34329196"
$wb.Worksheets.Item(3).Name = "summ35746218"
$wb.Worksheets.Item(4).Name = "summ36654416"
$wb.Worksheets.Item(5).Name = "summ37665193"
$wb.Worksheets.Item(6).Name = "summ38587992"
$wb.Worksheets.Item(7).Name = "summ39490908"
$wb.Worksheets.Item(8).Name = "summ41144416"
$wb.Worksheets.Item(9).Name = "summ41890856"

# Update coefficient table values on each sheet

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1522.956391481482
$ws.Cells.Item(2, 3).Value = 0.3139952371443505
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 320.8823705525579
$ws.Cells.Item(3, 3).Value = 0.7524474837196825
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 458.7944397425745
$ws.Cells.Item(4, 3).Value = 0.6387021542210245
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1394.065602666386
$ws.Cells.Item(5, 3).Value = 0.1522754643830427
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -16.95717021966171
$ws.Cells.Item(6, 3).Value = 0.6713623279424463
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1168.301481398934
$ws.Cells.Item(7, 3).Value = [double]"1.36118546263149e-25"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -12.31713694183266
$ws.Cells.Item(8, 3).Value = 0.01540983012917408
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 196.5557733039777
$ws.Cells.Item(9, 3).Value = [double]"6.683309498759672e-10"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 458.6530450215633
$ws.Cells.Item(10, 3).Value = [double]"1.264420813314676e-78"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01602808987721869
$ws.Cells.Item(11, 3).Value = 0.1989840766675158
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"1.748260547686058e-05"
$ws.Cells.Item(12, 3).Value = 0.3800868017112797
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -3.638924868144806
$ws.Cells.Item(13, 3).Value = 0.535867850291717
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 0.2874809862135845
$ws.Cells.Item(14, 3).Value = 0.970341550358405
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -824.029015429797
$ws.Cells.Item(15, 3).Value = 0.2670605582316218
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1446.093444853196
$ws.Cells.Item(16, 3).Value = 0.01009036588594601
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2610.827298122844
$ws.Cells.Item(17, 3).Value = [double]"8.631782862826166e-05"

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1404.816031785901
$ws.Cells.Item(2, 3).Value = 0.3538028766396497
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1126.412146346317
$ws.Cells.Item(3, 3).Value = 0.2606371419115787
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 957.1424792290356
$ws.Cells.Item(4, 3).Value = 0.3209397572282168
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1918.772861000336
$ws.Cells.Item(5, 3).Value = 0.0459062996700773
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -2.552053722425594
$ws.Cells.Item(6, 3).Value = 0.9494802920276616
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1261.955838750424
$ws.Cells.Item(7, 3).Value = [double]"2.280453896793739e-29"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -10.21096460672584
$ws.Cells.Item(8, 3).Value = 0.04441617351310496
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 210.2561839581641
$ws.Cells.Item(9, 3).Value = [double]"4.255875684020215e-11"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 437.3074963335171
$ws.Cells.Item(10, 3).Value = [double]"7.682177069536356e-71"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01321552489993552
$ws.Cells.Item(11, 3).Value = 0.2886396703033438
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"-1.237737283094702e-07"
$ws.Cells.Item(12, 3).Value = 0.9949997909967985
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -0.9848800188664431
$ws.Cells.Item(13, 3).Value = 0.8679904756999404
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -2.904889717555478
$ws.Cells.Item(14, 3).Value = 0.7129315531405407
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -204.1089602036714
$ws.Cells.Item(15, 3).Value = 0.7849797953802943
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1025.68765365046
$ws.Cells.Item(16, 3).Value = 0.0687244679533532
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2479.891094866841
$ws.Cells.Item(17, 3).Value = 0.0001950563371773189

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1100.523187334641
$ws.Cells.Item(2, 3).Value = 0.4799838085099545
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1579.186379370582
$ws.Cells.Item(3, 3).Value = 0.1358189852818083
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 1336.521727650496
$ws.Cells.Item(4, 3).Value = 0.1909702787153305
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 2261.069607196446
$ws.Cells.Item(5, 3).Value = 0.02645126220283478
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -25.72601136887959
$ws.Cells.Item(6, 3).Value = 0.5247628747394614
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1286.291355586122
$ws.Cells.Item(7, 3).Value = [double]"4.038282130341576e-30"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -14.92935871470974
$ws.Cells.Item(8, 3).Value = 0.003361563536107021
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 195.1088244188862
$ws.Cells.Item(9, 3).Value = [double]"7.784343387039226e-10"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 445.3611713462915
$ws.Cells.Item(10, 3).Value = [double]"5.749743858383188e-74"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01525459893861003
$ws.Cells.Item(11, 3).Value = 0.2241690487716188
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"1.239654988448021e-05"
$ws.Cells.Item(12, 3).Value = 0.5325842731313941
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 0.3412504079996624
$ws.Cells.Item(13, 3).Value = 0.9546869828587694
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 0.003900885306129087
$ws.Cells.Item(14, 3).Value = 0.9996072856717817
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -839.6549294786881
$ws.Cells.Item(15, 3).Value = 0.2631583754348311
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 894.58162133478
$ws.Cells.Item(16, 3).Value = 0.1169350314394295
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 1990.594511223455
$ws.Cells.Item(17, 3).Value = 0.002924132210709429

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1920.017397731989
$ws.Cells.Item(2, 3).Value = 0.2110554333661575
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1087.104056701781
$ws.Cells.Item(3, 3).Value = 0.3029781986674615
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 853.2517788545372
$ws.Cells.Item(4, 3).Value = 0.4013621306907046
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1789.669672177831
$ws.Cells.Item(5, 3).Value = 0.07737283382532333
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -4.148980630510355
$ws.Cells.Item(6, 3).Value = 0.9172743946619148
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1280.926951788521
$ws.Cells.Item(7, 3).Value = [double]"3.203533960392958e-30"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -6.768143240973282
$ws.Cells.Item(8, 3).Value = 0.1824175051866774
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 209.5702590891539
$ws.Cells.Item(9, 3).Value = [double]"4.132436748973043e-11"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 439.3791423782725
$ws.Cells.Item(10, 3).Value = [double]"3.746785566532748e-73"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01650605427805925
$ws.Cells.Item(11, 3).Value = 0.180450393673233
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"3.335018071883189e-06"
$ws.Cells.Item(12, 3).Value = 0.8646490566289007
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -6.133427030858371
$ws.Cells.Item(13, 3).Value = 0.2989898955677497
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -6.789744254166202
$ws.Cells.Item(14, 3).Value = 0.3885537969449007
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -375.4779977765078
$ws.Cells.Item(15, 3).Value = 0.6176850838492882
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1504.465673301621
$ws.Cells.Item(16, 3).Value = 0.007635467102107696
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2526.408274836029
$ws.Cells.Item(17, 3).Value = 0.000150852941912307

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = -583.7607534865236
$ws.Cells.Item(2, 3).Value = 0.7009291638821317
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 701.8199803292456
$ws.Cells.Item(3, 3).Value = 0.490432710561726
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 720.9027341693759
$ws.Cells.Item(4, 3).Value = 0.4613651905234
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1648.104257331252
$ws.Cells.Item(5, 3).Value = 0.09107016010562252
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = 33.26785073823717
$ws.Cells.Item(6, 3).Value = 0.4077329169443252
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1283.191500965684
$ws.Cells.Item(7, 3).Value = [double]"2.048786412906795e-30"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -9.983539793725027
$ws.Cells.Item(8, 3).Value = 0.04903987942922815
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 193.5213370067214
$ws.Cells.Item(9, 3).Value = [double]"1.548722125698586e-09"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 463.1685857756894
$ws.Cells.Item(10, 3).Value = [double]"1.483777548961002e-79"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.00738619292760432
$ws.Cells.Item(11, 3).Value = 0.5513678813143341
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"6.937432287975203e-06"
$ws.Cells.Item(12, 3).Value = 0.7256496529881957
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 1.685607307123138
$ws.Cells.Item(13, 3).Value = 0.7746113663561709
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 14.99864556718244
$ws.Cells.Item(14, 3).Value = 0.05487327246574399
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -831.233673112442
$ws.Cells.Item(15, 3).Value = 0.2646403412301429
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 842.6949332037902
$ws.Cells.Item(16, 3).Value = 0.1343658078268923
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 3019.448372924573
$ws.Cells.Item(17, 3).Value = [double]"5.072564459693287e-06"

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1321.551012888364
$ws.Cells.Item(2, 3).Value = 0.3781304087098454
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 825.1467438880245
$ws.Cells.Item(3, 3).Value = 0.4023210696987405
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 811.8835678358482
$ws.Cells.Item(4, 3).Value = 0.3905965149576786
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1772.181230239138
$ws.Cells.Item(5, 3).Value = 0.06002825330167061
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = 34.14671483028839
$ws.Cells.Item(6, 3).Value = 0.3914541204915185
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1301.302132492614
$ws.Cells.Item(7, 3).Value = [double]"2.762009776219429e-31"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -12.38260949759528
$ws.Cells.Item(8, 3).Value = 0.01510728087238344
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 203.9690648039056
$ws.Cells.Item(9, 3).Value = [double]"1.338962315153566e-10"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 471.9419050269357
$ws.Cells.Item(10, 3).Value = [double]"2.204613405157004e-83"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.008621667690382456
$ws.Cells.Item(11, 3).Value = 0.4863961067834947
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"1.355995381424888e-05"
$ws.Cells.Item(12, 3).Value = 0.4905412511695947
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -1.90246595912627
$ws.Cells.Item(13, 3).Value = 0.7468263946076299
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -1.725513361254613
$ws.Cells.Item(14, 3).Value = 0.8255463244984131
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -1127.411208918019
$ws.Cells.Item(15, 3).Value = 0.1328545196869658
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 713.3732558611631
$ws.Cells.Item(16, 3).Value = 0.2043936673964886
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2885.799778093011
$ws.Cells.Item(17, 3).Value = [double]"1.229466320156847e-05"

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 2378.744629915376
$ws.Cells.Item(2, 3).Value = 0.1255765683702301
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1430.452865394398
$ws.Cells.Item(3, 3).Value = 0.1754213353309697
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 1179.904706905374
$ws.Cells.Item(4, 3).Value = 0.2461705891027984
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 2061.464063880255
$ws.Cells.Item(5, 3).Value = 0.04201433883015469
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = 35.79310989491601
$ws.Cells.Item(6, 3).Value = 0.3731736831193847
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1252.632154323406
$ws.Cells.Item(7, 3).Value = [double]"6.15373871130937e-29"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -13.37350424632936
$ws.Cells.Item(8, 3).Value = 0.008527858060382665
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 179.6987092003609
$ws.Cells.Item(9, 3).Value = [double]"1.948864286292349e-08"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 429.8206527791115
$ws.Cells.Item(10, 3).Value = [double]"9.603540728143535e-69"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01206340081669634
$ws.Cells.Item(11, 3).Value = 0.3338851579192457
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"1.884267322278682e-06"
$ws.Cells.Item(12, 3).Value = 0.9240111187876345
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -9.476829565768709
$ws.Cells.Item(13, 3).Value = 0.1097681927061229
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -8.07743521248636
$ws.Cells.Item(14, 3).Value = 0.3050147701031933
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -306.4431397219337
$ws.Cells.Item(15, 3).Value = 0.6827376828280968
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1264.950835487229
$ws.Cells.Item(16, 3).Value = 0.02523635938138626
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2475.481153840133
$ws.Cells.Item(17, 3).Value = 0.0002036313543577991

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 842.65998369958
$ws.Cells.Item(2, 3).Value = 0.5780513488636563
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 968.5100474413645
$ws.Cells.Item(3, 3).Value = 0.3366265020367183
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 934.0906637930194
$ws.Cells.Item(4, 3).Value = 0.3349163451178361
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1851.429328487838
$ws.Cells.Item(5, 3).Value = 0.05508805762567355
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -9.447902150029858
$ws.Cells.Item(6, 3).Value = 0.8141476449883986
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1292.147883136498
$ws.Cells.Item(7, 3).Value = [double]"1.801871405756778e-30"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -9.141900177055616
$ws.Cells.Item(8, 3).Value = 0.07566016599418232
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 219.6929951021695
$ws.Cells.Item(9, 3).Value = [double]"7.011428462454314e-12"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 440.512025330628
$ws.Cells.Item(10, 3).Value = [double]"9.21133423963321e-72"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.0148525235816674
$ws.Cells.Item(11, 3).Value = 0.2351155738881615
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"4.518561805631733e-06"
$ws.Cells.Item(12, 3).Value = 0.8218001054100158
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -1.350915399002004
$ws.Cells.Item(13, 3).Value = 0.8196906654710817
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 0.8448845579522986
$ws.Cells.Item(14, 3).Value = 0.9143130239905707
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = 23.81687335776749
$ws.Cells.Item(15, 3).Value = 0.9749109809294129
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1270.139530749772
$ws.Cells.Item(16, 3).Value = 0.02495115511589098
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2576.539054058287
$ws.Cells.Item(17, 3).Value = 0.0001134298708987113

$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 687.9777614048743
$ws.Cells.Item(2, 3).Value = 0.6602167046583898
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 762.4317653088815
$ws.Cells.Item(3, 3).Value = 0.4768005841753322
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 782.8696901580112
$ws.Cells.Item(4, 3).Value = 0.4499579152266726
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1672.381919064916
$ws.Cells.Item(5, 3).Value = 0.1057574387590278
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = 7.037361516903964
$ws.Cells.Item(6, 3).Value = 0.860065284351569
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -1345.889224854322
$ws.Cells.Item(7, 3).Value = [double]"4.71518379310548e-33"
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -12.87794998100712
$ws.Cells.Item(8, 3).Value = 0.01183778730345415
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 187.2850230372941
$ws.Cells.Item(9, 3).Value = [double]"4.883779955009357e-09"
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 452.5471795430918
$ws.Cells.Item(10, 3).Value = [double]"7.665984996493868e-76"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.01675839794845372
$ws.Cells.Item(11, 3).Value = 0.1816080712211121
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"8.006066002605822e-06"
$ws.Cells.Item(12, 3).Value = 0.6876256378518751
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 0.06168009899576143
$ws.Cells.Item(13, 3).Value = 0.99170312513193
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 6.327466347512345
$ws.Cells.Item(14, 3).Value = 0.4232243030666653
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -1300.696503429195
$ws.Cells.Item(15, 3).Value = 0.08556153773223489
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = 1409.02292054315
$ws.Cells.Item(16, 3).Value = 0.01280370843125633
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = 2772.029021112467
$ws.Cells.Item(17, 3).Value = [double]"3.21370152738139e-05"
